# Update cryptos list - GitHub Actions scheduled refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddress, $value) {
    # Force the cell to keep a plain text representation instead of letting
    # Excel auto-convert numeric-looking strings (e.g. "429.71") into a
    # floating point number, then restore the original (default) style.
    $rng = $ws.Range($rangeAddress)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "66.560.26"
$ws.Range("E2").Value = "  -0.49%  "

# Row 3 - Ethereum
Set-TextValue "D3" "3.849.97"
$ws.Range("E3").Value = "  +3.31%  "

# Row 4 - TetherUSD
Set-TextValue "D4" "0.999"
$ws.Range("E4").Value = "  -0.21%  "

# Row 5 - BNB
Set-TextValue "D5" "429.71"
$ws.Range("E5").Value = "  +2.09%  "

# Row 6 - Solana
Set-TextValue "D6" "130.93"
$ws.Range("E6").Value = "  -0.88%  "

# Row 7 - LidoStakedEther
Set-TextValue "D7" "3.843.71"
$ws.Range("E7").Value = "  +3.32%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  -5.39%  "

# Row 9 - USDC
Set-TextValue "D9" "0.999"
$ws.Range("E9").Value = "  -0.04%  "

# Row 10 - Cardano
$ws.Range("E10").Value = "  -6.14%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value = "  -8.63%  "

# Row 12 - ShibaInu
Set-TextValue "D12" "0.0000364"
$ws.Range("E12").Value = "  -10.56%  "

# Row 13 - Avalanche
Set-TextValue "D13" "40.81"
$ws.Range("E13").Value = "  -5.24%  "

# Row 14 - WrappedliquidstakedEther2.0
Set-TextValue "D14" "4.444.74"
$ws.Range("E14").Value = "  +3.54%  "

# Row 15 - Polkadot
$ws.Range("E15").Value = "  -5.35%  "

# Row 16 - Uniswap
Set-TextValue "D16" "15.76"
$ws.Range("E16").Value = "  +18.43%  "

# Row 17 & 18 - TRON and WrappedEther swap places (WrappedEther now row 17, TRON row 18)
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D17" "3.858.29"
$ws.Range("E17").Value = "  +3.66%  "

$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue "D18" "0.138"
$ws.Range("E18").Value = "  -1.40%  "

# Row 19 - Chainlink
Set-TextValue "D19" "19.58"
$ws.Range("E19").Value = "  -5.87%  "

# Row 20 - WrappedBTC
Set-TextValue "D20" "66.961.71"
$ws.Range("E20").Value = "  +0.04%  "

# Row 21 - Polygon
$ws.Range("E21").Value = "  -6.73%  "

# Row 22 - BitcoinCash
Set-TextValue "D22" "408.31"
$ws.Range("E22").Value = "  -9.15%  "

# Row 23 - InternetComputer(DFINITY)
$ws.Range("E23").Value = "  -13.51%  "

# Row 24 - Litecoin
Set-TextValue "D24" "85.21"
$ws.Range("E24").Value = "  -5.43%  "

# Row 25 - ImmutableX
Set-TextValue "D25" "3.04"
$ws.Range("E25").Value = "  -4.08%  "

# Row 26 - EthereumClassic
Set-TextValue "D26" "37.10"
$ws.Range("E26").Value = "  -2.63%  "

# Row 27 - LEO
$ws.Range("E27").Value = "  +12.33%  "

# Row 28 - PancakeSwap
Set-TextValue "D28" "3.25"
$ws.Range("E28").Value = "  -2.80%  "

# Row 29 - Filecoin
Set-TextValue "D29" "9.49"
$ws.Range("E29").Value = "  -7.02%  "

# Row 30 - Bittensor
Set-TextValue "D30" "688.16"
$ws.Range("E30").Value = "  +4.81%  "

# Row 31 - Cosmos
Set-TextValue "D31" "12.46"
$ws.Range("E31").Value = "  -2.62%  "

# Row 32 - Hedera
$ws.Range("E32").Value = "  -3.02%  "

# Row 33 - Toncoin
$ws.Range("E33").Value = "  -3.29%  "

# Row 34 - RenderToken
Set-TextValue "D34" "7.15"
$ws.Range("E34").Value = "  -1.97%  "

# Row 35 - Kaspa
Set-TextValue "D35" "0.152"
$ws.Range("E35").Value = "  -8.22%  "

# Row 36 - InjectiveProtocol
Set-TextValue "D36" "38.72"
$ws.Range("E36").Value = "  -7.89%  "

# Row 37 - Dai
Set-TextValue "D37" "1.00"
$ws.Range("E37").Value = "  +0.04%  "

# Row 38 - PEPE
Set-TextValue "D38" "0.0₃0798"
$ws.Range("E38").Value = "  +6.17%  "

# Row 39 - OKB
Set-TextValue "D39" "55.19"
$ws.Range("E39").Value = "  -3.67%  "

# Row 40 - ThetaToken
Set-TextValue "D40" "3.09"
$ws.Range("E40").Value = "  -0.39%  "

# Row 41 - VeChain
$ws.Range("E41").Value = "  -8.21%  "

# Row 42 - FirstDigitalUSD
$ws.Range("E42").Value = "  +0.40%  "

# Row 43 & 44 - Stellar and NEARProtocol swap places (NEARProtocol now row 43, Stellar row 44)
$ws.Range("B43").Value = "NEARProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue "D43" "4.55"
$ws.Range("E43").Value = "  +3.22%  "

$ws.Range("B44").Value = "Stellar"
$ws.Range("C44").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue "D44" "0.137"
$ws.Range("E44").Value = "  -8.54%  "

# Row 45 - Monero (price only)
Set-TextValue "D45" "148.35"

# Row 46 - LidoDAOToken
Set-TextValue "D46" "3.30"
$ws.Range("E46").Value = "  -4.64%  "

# Row 47 & 48 - ApeXProtocol and ARBITRUM swap places (ARBITRUM now row 47, ApeXProtocol row 48)
$ws.Range("B47").Value = "ARBITRUM"
$ws.Range("C47").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue "D47" "2.08"
$ws.Range("E47").Value = "  -2.71%  "

$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
Set-TextValue "D48" "3.11"
$ws.Range("E48").Value = "  -5.54%  "

# Row 49 - EnergySwap
Set-TextValue "D49" "26.22"
$ws.Range("E49").Value = "  -8.82%  "

# Row 50 - Stacks
Set-TextValue "D50" "2.79"
$ws.Range("E50").Value = "  -4.17%  "

# Row 51 - WEMIXToken
Set-TextValue "D51" "2.53"
$ws.Range("E51").Value = "  -5.12%  "
